# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-07 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-08 Monday", 2)

# Update the division problems in the table, addressed by (row, column)
# so that values which coincidentally collide with other old/new values
# (e.g. "83÷4=" is both an old value at (17,1) and the new value at (1,3))
# are never ambiguous.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "87÷4="
$t.Cell(1,2).Range.Text  = "55÷5="
$t.Cell(1,3).Range.Text  = "83÷4="
$t.Cell(1,4).Range.Text  = "85÷3="
$t.Cell(1,5).Range.Text  = "70÷8="

$t.Cell(5,1).Range.Text  = "89÷8="
$t.Cell(5,2).Range.Text  = "70÷6="
$t.Cell(5,3).Range.Text  = "74÷6="
$t.Cell(5,4).Range.Text  = "61÷8="
$t.Cell(5,5).Range.Text  = "41÷8="

$t.Cell(9,1).Range.Text  = "17÷7="
$t.Cell(9,2).Range.Text  = "31÷6="
$t.Cell(9,3).Range.Text  = "58÷6="
$t.Cell(9,4).Range.Text  = "40÷2="
$t.Cell(9,5).Range.Text  = "57÷6="

$t.Cell(13,1).Range.Text = "68÷9="
$t.Cell(13,2).Range.Text = "54÷6="
$t.Cell(13,3).Range.Text = "64÷2="
$t.Cell(13,4).Range.Text = "98÷2="
$t.Cell(13,5).Range.Text = "75÷4="

$t.Cell(17,1).Range.Text = "15÷7="
$t.Cell(17,2).Range.Text = "46÷2="
$t.Cell(17,3).Range.Text = "54÷7="
$t.Cell(17,4).Range.Text = "12÷8="
$t.Cell(17,5).Range.Text = "94÷6="
